$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Users")

# Insert a new row at position 4, shifting old row 4 (F00043/043) down to row 5
$ws.Rows.Item(4).Insert()

# Copy the text format (style index 4 equivalent: numFmt "@" + right aligned)
# from an existing formatted cell (C3) onto the new C-column cells so that the
# numeric-looking values (073, 074, 168, 037) are preserved as text, reusing
# the existing style instead of creating a new one.
$ws.Range("C3").Copy()
$ws.Range("C4").PasteSpecial(-4122)
$ws.Range("C6:C12").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Column A values first (matches the order new shared strings were authored in)
$ws.Cells.Item(4, 1).Value = "F00289"
$ws.Cells.Item(6, 1).Value = "F00073"
$ws.Cells.Item(7, 1).Value = "F00473"
$ws.Cells.Item(8, 1).Value = "F00474"
$ws.Cells.Item(9, 1).Value = "ANOVELLO"
$ws.Cells.Item(10, 1).Value = "F00274"

# Column C values for rows 4-10
$ws.Cells.Item(4, 3).Value = "089"
$ws.Cells.Item(6, 3).Value = "073"
$ws.Cells.Item(7, 3).Value = "073"
$ws.Cells.Item(8, 3).Value = "074"
$ws.Cells.Item(9, 3).Value = "Gerente Operativo"
$ws.Cells.Item(10, 3).Value = "074"

# Rows 11-12: column A then column C
$ws.Cells.Item(11, 1).Value = "F02653"
$ws.Cells.Item(12, 1).Value = "F00743"
$ws.Cells.Item(11, 3).Value = "168"
$ws.Cells.Item(12, 3).Value = "037"

# Column C width (stored width of 20 character-units; ColumnWidth/stored-width
# have a small fixed offset in this runtime, so 19.2 round-trips to width="20")
$ws.Columns.Item(3).ColumnWidth = 19.2

# Row 2 custom height
$ws.Rows.Item(2).RowHeight = 14.25

# Selection on Users sheet
$ws.Range("D14").Select()

# Make Users the active/selected tab (was Modulos before)
$ws.Select()

Write-Host "Done"
